# crsr shares / covered calls
$wb = $excel.ActiveWorkbook

$wsShares = $wb.Worksheets.Item("CRSRShares")
$wsCalls  = $wb.Worksheets.Item("CRSRCoveredCalls")

# --- CRSRShares: append a new holding row (row 14) ---
$wsShares.Cells.Item(14, 1).Value = 1
$wsShares.Cells.Item(14, 2).Value = 29.62

# --- CRSRCoveredCalls: premium paid on row 10 updated ---
$wsCalls.Cells.Item(10, 2).Value = 27

# --- Update selections on each sheet ---
[void]$wsShares.Range("B15").Select()
[void]$wsCalls.Range("B11").Select()

# --- CRSRShares is now the active/visible tab ---
[void]$wsShares.Activate()
